$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.900.16"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "2.262.95"
$ws.Range("E3").Value = "  +2.18%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "270.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.97%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0933"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.50%  "

$ws.Range("E13").Value = "  +0.75%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.603.25"
$ws.Range("E14").Value = "  +2.24%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.826"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.30%  "

$ws.Range("D17").Value = "2.243.55"
$ws.Range("E17").Value = "  +1.64%  "

$ws.Range("D18").Value = "43.855.56"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("E19").Value = "  +1.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.59%  "

$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.73%  "

$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.18%  "

$ws.Range("E27").Value = "  +10.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("E29").Value = "  -5.02%  "

$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("E33").Value = "  +4.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0350"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.69%  "

$ws.Range("E39").Value = "  +19.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.250"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.33%  "

$ws.Range("E45").Value = "  +4.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.00%  "

$ws.Range("E47").Value = "  +7.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "99.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.428"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.71%  "

$ws.Range("D51").Value = "2.484.98"
$ws.Range("E51").Value = "  +2.20%  "
